# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (commit: "Updated cryptos list ... with GitHub Actions").
# - Refreshes Price (column D) and Volume(1h) (column E) figures for most rows.
# - Re-sorts four rows (22-25) whose relative ranking changed, by rewriting
#   Coin/Link/Price/Volume for those rows in the new order:
#     22: InternetComputer(DFINITY), 23: Litecoin, 24: PancakeSwap, 25: Toncoin
#
# Note: several Price values are plain decimal numbers (e.g. "0.999"); a
# leading apostrophe is used so Excel stores them as text (matching the
# original column's text formatting) instead of auto-converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.334.09"
$ws.Range("E2").Value = "  -2.08%  "

$ws.Range("D3").Value = "3.478.52"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'611.47"
$ws.Range("E5").Value = "  +5.17%  "

$ws.Range("D6").Value = "'189.06"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  -0.55%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  -3.91%  "

$ws.Range("D10").Value = "'0.647"
$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("D11").Value = "'52.74"

$ws.Range("E12").Value = "  -3.93%  "

$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").Value = "4.039.44"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("D15").Value = "'611.25"
$ws.Range("E15").Value = "  +6.18%  "

$ws.Range("D16").Value = "69.446.44"
$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("D18").Value = "'12.49"
$ws.Range("E18").Value = "  -2.94%  "

$ws.Range("D19").Value = "3.485.42"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("E21").Value = "  -2.43%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'17.08"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'106.37"
$ws.Range("E23").Value = "  +12.94%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.69"
$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'5.12"
$ws.Range("E25").Value = "  +4.71%  "

$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  -2.95%  "

$ws.Range("D28").Value = "'9.65"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("D29").Value = "'33.33"
$ws.Range("E29").Value = "  +1.57%  "

$ws.Range("D30").Value = "'6.89"
$ws.Range("E30").Value = "  -4.84%  "

$ws.Range("D31").Value = "'12.53"
$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("D32").Value = "'3.93"
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").Value = "'63.13"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("E35").Value = "  -5.86%  "

$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").Value = "3.612.18"
$ws.Range("E37").Value = "  -0.55%  "

$ws.Range("D38").Value = "'3.62"
$ws.Range("E38").Value = "  +4.80%  "

$ws.Range("E39").Value = "  -4.86%  "

$ws.Range("D40").Value = "'504.10"
$ws.Range("E40").Value = "  -6.48%  "

$ws.Range("D41").Value = "'36.47"
$ws.Range("E41").Value = "  -4.51%  "

$ws.Range("E42").Value = "  -6.55%  "

$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  -2.19%  "

$ws.Range("D46").Value = "'0.141"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  -4.35%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").Value = "'8.68"
$ws.Range("E49").Value = "  -7.05%  "

$ws.Range("D50").Value = "'130.86"
$ws.Range("E50").Value = "  -4.27%  "

$ws.Range("E51").Value = "  -7.50%  "

